$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$headers = @("id", "street", "postcode", "X", "Y")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ids = @(1, 2, 3, 4)
$streets = @("High Street", "Low Street", "Back Lane", "Top Drive")
$postcodes = @("A99 1AA", "B1 2BD", "C4 8FG", "D5 7YG")
$xs = @(525141, 392276, 425492, 261626)
$ys = @(365969, 214282, 430187, 114633)

for ($r = 0; $r -lt $ids.Length; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = $ids[$r]
}
for ($r = 0; $r -lt $streets.Length; $r++) {
    $ws.Cells.Item($r + 2, 2).Value = $streets[$r]
}
for ($r = 0; $r -lt $postcodes.Length; $r++) {
    $ws.Cells.Item($r + 2, 3).Value = $postcodes[$r]
}
for ($r = 0; $r -lt $xs.Length; $r++) {
    $ws.Cells.Item($r + 2, 4).Value = $xs[$r]
}
for ($r = 0; $r -lt $ys.Length; $r++) {
    $ws.Cells.Item($r + 2, 5).Value = $ys[$r]
}

$ws.Range("C6").Select()
